$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TeamStats")

# Add the missing AVERAGE formula for the velocity column (matches the
# pattern already used by the storyPoints/duration/idlePeriod columns).
$ws.Range("B17").Formula = "=AVERAGE(B2:B15)"

# Update the active selection on the sheet to match the latest edit:
# range B17:C17 selected (active cell lands on the range's top-left, B17).
$ws.Activate()
$ws.Range("B17:C17").Select()
